$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1395
$ws.Range("I2").Value = 3824
$ws.Range("J2").Value = 15512
$ws.Range("K2").Value = 82
$ws.Range("L2").Value = 4327
$ws.Range("M2").Value = 283
$ws.Range("N2").Value = 2684
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 60
$ws.Range("Q2").Value = 26
$ws.Range("R2").Value = 203
$ws.Range("S2").Value = 1703
$ws.Range("T2").Value = 2749
$ws.Range("U2").Value = 195
$ws.Range("V2").Value = 24247
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 23948
$ws.Range("Y2").Value = 42
$ws.Range("Z2").Value = 338
$ws.Range("AA2").Value = 161
